$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-6 (shifts remaining data rows up by 5)
$ws.Range("A2:C6").EntireRow.Delete() | Out-Null

# Add new rows of data at the end (now rows 18-21, since old row22 data is now at row17)
$newData = @(
    @(-0.2174680233001709, 1.732872009277344, -2.470797300338745),
    @(1.422094345092773, 4.435187339782715, 0.8324565887451172),
    @(2.657721996307373, 5.389358520507812, 2.652835130691528),
    @(0.5099197626113892, 0.4361577928066253, 2.850908041000366)
)

$startRow = 18
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
